$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 567, shifting existing rows 567-616 down to 569-618.
$ws.Rows.Item(567).Resize(2).Insert()

# Copy formatting (style) of column D date cells down so the new rows inherit the date style.
$ws.Range("D569").Copy()
$ws.Range("D567:D568").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 567: new record
$ws.Range("A567").Value = 3
$ws.Range("B567").Value = "Femacal de La Calera"
$ws.Range("C567").Value = "Coquimbo"
$ws.Range("D567").Value = 44578
$ws.Range("E567").Value = 5
$ws.Range("F567").Value = "Fruta"
$ws.Range("G567").Value = 100102
$ws.Range("H567").Value = "Cítricos"
$ws.Range("I567").Value = 100102005
$ws.Range("J567").Value = "Naranja"
$ws.Range("K567").Value = "Valencia"
$ws.Range("L567").Value = "Primera"
$ws.Range("M567").Value = 80
$ws.Range("N567").Value = 6500
$ws.Range("O567").Value = 6500
$ws.Range("P567").Value = 6500
$ws.Range("Q567").Value = "`$/malla 13 kilos"
$ws.Range("R567").Value = "Provincia de Quillota"
$ws.Range("S567").Value = 500
$ws.Range("T567").Value = 13

# Row 568: new record
$ws.Range("A568").Value = 3
$ws.Range("B568").Value = "Femacal de La Calera"
$ws.Range("C568").Value = "Coquimbo"
$ws.Range("D568").Value = 44578
$ws.Range("E568").Value = 5
$ws.Range("F568").Value = "Fruta"
$ws.Range("G568").Value = 100102
$ws.Range("H568").Value = "Cítricos"
$ws.Range("I568").Value = 100102005
$ws.Range("J568").Value = "Naranja"
$ws.Range("K568").Value = "Valencia"
$ws.Range("L568").Value = "Segunda"
$ws.Range("M568").Value = 90
$ws.Range("N568").Value = 5000
$ws.Range("O568").Value = 5000
$ws.Range("P568").Value = 5000
$ws.Range("Q568").Value = "`$/malla 13 kilos"
$ws.Range("R568").Value = "Provincia de Quillota"
$ws.Range("S568").Value = 385
$ws.Range("T568").Value = 13
